$d = $word.ActiveDocument
$sec = $d.Sections(1)

# --- First-page header: BTec_Logo-Orange picture ---
# wp:docPr / pic:cNvPr name: image2.jpg -> image1.jpg
$hdr = $sec.Headers(2)
if ($hdr.Exists -and ($hdr.Range.InlineShapes.Count -ge 1)) {
    $hdrShape = $hdr.Range.InlineShapes(1)
    $hdrShape.Name = "image1.jpg"
}

# --- Footers: PearsonLogo picture (appears in both the primary and the
#     first-page footer) ---
# wp:docPr / pic:cNvPr name: image1.png -> image2.png
for ($i = 1; $i -le 3; $i++) {
    $ftr = $sec.Footers($i)
    if ($ftr.Exists -and ($ftr.Range.InlineShapes.Count -ge 1)) {
        $ftrShape = $ftr.Range.InlineShapes(1)
        $ftrShape.Name = "image2.png"
    }
}
